# Update to Devices in Dashboard
# Adds a new "DTYPE" column to the Devices table (Table1), classifying each
# row as Server / Node / Unregistered, and switches the active sheet/tab
# selection from Alerts back to Devices.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Devices")
$tbl = $ws.ListObjects.Item("Table1")

# --- Add the new table column (becomes column L / id 12) ------------------
$newCol = $tbl.ListColumns.Add()

# Header
$ws.Range("L1").Value = "DTYPE"

# Give the new column an explicit width (matches how Excel records a width
# once a previously "default width" column has been touched).
$ws.Columns.Item(12).ColumnWidth = 8.14

# --- Fill in the per-row classification ------------------------------------
$serverRows = @(2, 3)
$unregisteredRows = @(31, 32, 33)

for ($r = 2; $r -le 33; $r++) {
    if ($serverRows -contains $r) {
        $value = "Server"
    } elseif ($unregisteredRows -contains $r) {
        $value = "Unregistered"
    } else {
        $value = "Node"
    }
    $ws.Range("L$r").Value = $value
}

# --- Switch the active sheet/selection back to Devices ---------------------
[void]$ws.Activate()
[void]$ws.Range("L29").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
